$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Gillette Mach3 2 unidades
$ws.Range("B2").Value = "Carga para Aparelho de Barbear Gillette Mach3 com 2 unidades"
$ws.Range("C2").ClearContents()
$ws.Range("E2").Formula = '=HYPERLINK("https://www.mercadolivre.com.br/carga-para-lmina-de-barbear-gillette-mach3-2-unidades/p/MLB17355357", "11.71")'
$ws.Range("H2").Value = 11.71

# Row 3 - Gillette Fusion5 4 unidades
$ws.Range("B3").Value = "Carga para Aparelho de Barbear Gillette Fusion 5 - 4 unidades"
$ws.Range("C3").ClearContents()
$ws.Range("E3").Formula = '=HYPERLINK("https://www.mercadolivre.com.br/4-cartuchos-refil-para-aparelho-fusion-5-gillette/p/MLB16088319", "87.02")'
$ws.Range("H3").Value = 87.02

# Row 4 - Gillette Mach3 1 unidade
$ws.Range("B4").Value = "Aparelho de Barbear Gillette Mach3 com 1 Unidade"
$ws.Range("C4").ClearContents()
$ws.Range("E4").Formula = '=HYPERLINK("https://www.mercadolivre.com.br/aparelho-de-barbear-mach3-carbono-reutilizavel/p/MLB23207098", "26.59")'
$ws.Range("H4").Value = 26.59
